$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The localization status moved from "Ready for handoff" to "In Translation".
# That shared text shows up in the Overview rollup (one column per locale)
# and as the "Status" column on each per-locale report sheet.
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# With the shorter status text, those columns no longer need to be as wide.
$wsOverview.Columns.Item(5).ColumnWidth = 13.4101845877511
$wsOverview.Columns.Item(6).ColumnWidth = 13.4101845877511
$wsZhCn.Columns.Item(3).ColumnWidth = 13.4101845877511
$wsDeDe.Columns.Item(3).ColumnWidth = 13.4101845877511
